$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "https://www.city.takamatsu.kagawa.jp/kurashi/kosodate/bunka/rekishi/index.html"
$ws.Range("G4").Value = "https://www.city.takamatsu.kagawa.jp/kurashi/kosodate/bunka/ishimin/index.html"
$ws.Range("G5").Value = "https://www.city.takamatsu.kagawa.jp/kurashi/kosodate/bunka/kounanrekishi/index.html"
$ws.Range("G6").Value = "https://www.city.takamatsu.kagawa.jp/kurashi/kosodate/bunka/sanuki.html"
